# Update version number from 0.1 to 1.0.
# A direct $ws.Range("D2").Value = "1.0" would make Excel auto-coerce the
# numeric-looking text into a plain number (and flip the cell's stored
# type), so instead we write it as a text formula on a scratch cell, copy
# it, and paste-special just the resulting value back onto D2. That keeps
# the cell's original style/type (a text cell) intact.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$scratch = $ws.Range("Z100")
$scratch.Formula = "=""1.0"""
$scratch.Copy()
$ws.Range("D2").PasteSpecial(-4163)  # xlPasteValues
$scratch.ClearContents()

# Swap the "Expected Results" text between TC1's last step (D13) and
# TC3's last step (D32): the success message now belongs to TC1's flow
# and the failure message now belongs to TC3's flow.
$tc1Text = $ws.Range("D13").Text
$tc3Text = $ws.Range("D32").Text

$ws.Range("D13").Value = $tc3Text
$ws.Range("D32").Value = $tc1Text
